$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new tracker row for 25/08/2016
$ws.Range("A6").Value = "25/08/2016"
$ws.Range("H6").Value = "completed"
$ws.Range("I6").Value = "in progress"

# Move active selection to the last edited cell, matching the committed file
$ws.Range("I6").Select()
